$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.163.20"
$ws.Range("E2").Value = "  -0.86%  "

$ws.Range("D3").Value = "3.429.33"
$ws.Range("E3").Value = "  -1.06%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "408.03"
$ws.Range("D5").Style = $ws.Range("B5").Style
$ws.Range("E5").Value = "  -0.64%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.50"
$ws.Range("D6").Style = $ws.Range("B6").Style
$ws.Range("E6").Value = "  +3.60%  "

$ws.Range("E7").Value = "  -0.18%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").Style = $ws.Range("B8").Style
$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("E9").Value = "  -2.26%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.122"
$ws.Range("D10").Style = $ws.Range("B10").Style
$ws.Range("E10").Value = "  -4.65%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.50"
$ws.Range("D11").Style = $ws.Range("B11").Style
$ws.Range("E11").Value = "  -2.68%  "

$ws.Range("E12").Value = "  -1.50%  "

$ws.Range("D13").Value = "3.963.56"
$ws.Range("E13").Value = "  -1.15%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.43"
$ws.Range("D14").Style = $ws.Range("B14").Style
$ws.Range("E14").Value = "  -3.38%  "

$ws.Range("E15").Value = "  -1.46%  "

$ws.Range("D16").Value = "3.419.93"
$ws.Range("E16").Value = "  -2.33%  "

$ws.Range("D17").Value = "62.146.29"
$ws.Range("E17").Value = "  -0.65%  "

$ws.Range("E18").Value = "  -3.10%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.03"
$ws.Range("D19").Style = $ws.Range("B19").Style
$ws.Range("E19").Value = "  +0.76%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000131"
$ws.Range("D20").Style = $ws.Range("B20").Style
$ws.Range("E20").Value = "  -4.05%  "

$ws.Range("E21").Value = "  -4.78%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "84.85"
$ws.Range("D22").Style = $ws.Range("B22").Style
$ws.Range("E22").Value = "  +3.35%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "315.41"
$ws.Range("D23").Style = $ws.Range("B23").Style
$ws.Range("E23").Value = "  +0.91%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.83"
$ws.Range("D24").Style = $ws.Range("B24").Style
$ws.Range("E24").Value = "  -2.97%  "

$ws.Range("E25").Value = "  -2.17%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.79"
$ws.Range("D26").Style = $ws.Range("B26").Style
$ws.Range("E26").Value = "  +9.70%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "29.75"
$ws.Range("D27").Style = $ws.Range("B27").Style
$ws.Range("E27").Value = "  -1.97%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.24"
$ws.Range("D28").Style = $ws.Range("B28").Style
$ws.Range("E28").Value = "  +1.80%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.71"
$ws.Range("D29").Style = $ws.Range("B29").Style
$ws.Range("E29").Value = "  -1.79%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.76"
$ws.Range("D30").Style = $ws.Range("B30").Style
$ws.Range("E30").Value = "  +2.43%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.175"
$ws.Range("D31").Style = $ws.Range("B31").Style
$ws.Range("E31").Value = "  -2.02%  "

$ws.Range("E32").Value = "  -4.33%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "42.91"
$ws.Range("D33").Style = $ws.Range("B33").Style
$ws.Range("E33").Value = "  -4.59%  "

$ws.Range("E34").Value = "  -0.12%  "

$ws.Range("E35").Value = "  -4.95%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0485"
$ws.Range("D36").Style = $ws.Range("B36").Style
$ws.Range("E36").Value = "  -1.53%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "52.20"
$ws.Range("D37").Style = $ws.Range("B37").Style
$ws.Range("E37").Value = "  -0.97%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.998"
$ws.Range("D38").Style = $ws.Range("B38").Style
$ws.Range("E38").Value = "  +0.17%  "

$ws.Range("E40").Value = "  -1.08%  "

$ws.Range("E41").Value = "  -0.25%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "137.79"
$ws.Range("D42").Style = $ws.Range("B42").Style
$ws.Range("E42").Value = "  -0.03%  "

$ws.Range("E43").Value = "  -0.19%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.294"
$ws.Range("D44").Style = $ws.Range("B44").Style
$ws.Range("E44").Value = "  +0.85%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.00"
$ws.Range("D45").Style = $ws.Range("B45").Style
$ws.Range("E45").Value = "  +0.25%  "

$ws.Range("E46").Value = "  -7.41%  "

$ws.Range("E47").Value = "  -1.67%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "21.45"
$ws.Range("D48").Style = $ws.Range("B48").Style
$ws.Range("E48").Value = "  -5.24%  "

$ws.Range("D49").Value = "2.131.96"
$ws.Range("E49").Value = "  -5.00%  "

$ws.Range("E50").Value = "  -3.83%  "

$ws.Range("E51").Value = "  +1.82%  "
